$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "AmbiguousAgeRange" / "AmbiguousAgeRangeObs" columns (N:O).
# This shifts every later section two columns to the left:
#   P(SurvInRep)->N, R(MeanAgeFirstRep)->P, T(StageFirstRep)->R,
#   V(ReproWithMaturation)->T, X(LongStages)->V, Z(MeanStageDuration)->X,
#   AB(VarStateDuration)->Z, AD(GrowthTransition)->AB, AE->AC
$ws.Range("N1:O1").EntireColumn.Delete()

# New reference-list entries for the (now relocated) ReproWithMaturation column.
$ws.Range("T2").Value = "Yes"
$ws.Range("T3").Value = "No"
$ws.Range("T4").Value = "Unknown"

# CensusType reference list gains an "Ambiguous" option (inserted before "Unknown").
$ws.Range("L7").Value = "Ambiguous"
$ws.Range("L8").Value = "Unknown"

# Reproduce the author's final selection/scroll position.
$ws.Range("T4").Select()
